# Scheduled-runner style refresh of market/profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N) across
# the per-class "Profits" tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2116.9
$ws.Range("I19").Value = 1587.25
$ws.Range("J19").Value = 2470
$ws.Range("K19").Value = 1587.25
$ws.Range("L19").Value = 2470
$ws.Range("M19").Value = -1412.25
$ws.Range("N19").Value = -2820

# ALC!row38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 312.92856
$ws.Range("I38").Value = 106.23077
$ws.Range("K38").Value = 318.69231
$ws.Range("M38").Value = 53.30768999999998

# ALC!row55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1634.5714
$ws.Range("I55").Value = 488
$ws.Range("J55").Value = 4501
$ws.Range("K55").Value = 488
$ws.Range("L55").Value = 4501
$ws.Range("M55").Value = -274
$ws.Range("N55").Value = -4929

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4683.6304
$ws.Range("I98").Value = 4852.25
$ws.Range("K98").Value = 4852.25
$ws.Range("M98").Value = -3354.25

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I113").Value = 111113030
$ws.Range("J113").Value = 3150
$ws.Range("K113").Value = 111113030
$ws.Range("L113").Value = 3150
$ws.Range("M113").Value = -111109776
$ws.Range("N113").Value = -9658

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4683.6304
$ws.Range("I122").Value = 4852.25
$ws.Range("K122").Value = 14556.75
$ws.Range("M122").Value = -12106.75

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4242.484
$ws.Range("I74").Value = 2358.842
$ws.Range("K74").Value = 2358.842
$ws.Range("M74").Value = -1484.842

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4242.484
$ws.Range("I77").Value = 2358.842
$ws.Range("K77").Value = 11794.21
$ws.Range("M77").Value = -7426.210000000001

# BSM!row38
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

# BSM!row62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 27499.5
$ws.Range("I62").Value = 27499
$ws.Range("J62").Value = 27500
$ws.Range("K62").Value = 27499
$ws.Range("L62").Value = 27500
$ws.Range("M62").Value = -26813
$ws.Range("N62").Value = -28872

# BSM!row65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 27499.5
$ws.Range("I65").Value = 27499
$ws.Range("J65").Value = 27500
$ws.Range("K65").Value = 82497
$ws.Range("L65").Value = 82500
$ws.Range("M65").Value = -79065
$ws.Range("N65").Value = -89364

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2693.1667
$ws.Range("I94").Value = 2176.077
$ws.Range("K94").Value = 2176.077
$ws.Range("M94").Value = -1725.077

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3982.7778
$ws.Range("I105").Value = 4028.8823
$ws.Range("K105").Value = 4028.8823
$ws.Range("M105").Value = -2281.8823

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2345.0908
$ws.Range("I107").Value = 1879.9062
$ws.Range("J107").Value = 3585.5833
$ws.Range("K107").Value = 1879.9062
$ws.Range("L107").Value = 3585.5833
$ws.Range("M107").Value = 40.0938000000001
$ws.Range("N107").Value = -7425.5833

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5547.615
$ws.Range("I31").Value = 4476.0586
$ws.Range("J31").Value = 7571.6665
$ws.Range("K31").Value = 4476.0586
$ws.Range("L31").Value = 7571.6665
$ws.Range("M31").Value = -4181.0586
$ws.Range("N31").Value = -8161.6665

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5547.615
$ws.Range("I34").Value = 4476.0586
$ws.Range("J34").Value = 7571.6665
$ws.Range("K34").Value = 4476.0586
$ws.Range("L34").Value = 7571.6665
$ws.Range("M34").Value = -4274.0586
$ws.Range("N34").Value = -7975.6665

# CRP!row43
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 24000
$ws.Range("J43").Value = 24000
$ws.Range("L43").Value = 24000
$ws.Range("N43").Value = -24368

# CRP!row53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8366.583000000001
$ws.Range("I58").Value = 3628.5715
$ws.Range("K58").Value = 3628.5715
$ws.Range("M58").Value = -3425.5715

# CRP!row101
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 24000
$ws.Range("J101").Value = 24000
$ws.Range("L101").Value = 24000
$ws.Range("N101").Value = -30490

# CRP!row111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 70640
$ws.Range("J111").Value = 70640
$ws.Range("L111").Value = 70640
$ws.Range("N111").Value = -78820

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8366.583000000001
$ws.Range("I136").Value = 3628.5715
$ws.Range("K136").Value = 10885.7145
$ws.Range("M136").Value = -8335.7145

# CRP!row137
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 84436
$ws.Range("J137").Value = 84436
$ws.Range("L137").Value = 84436
$ws.Range("N137").Value = -94636

# CUL!row14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 355.5
$ws.Range("I14").Value = 355.5
$ws.Range("K14").Value = 1066.5
$ws.Range("M14").Value = -893.5

# CUL!row50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 274.7143
$ws.Range("I50").Value = 62.5
$ws.Range("K50").Value = 187.5
$ws.Range("M50").Value = 293.5

# CUL!row53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 274.7143
$ws.Range("I53").Value = 62.5
$ws.Range("K53").Value = 187.5
$ws.Range("M53").Value = 293.5

# CUL!row81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 10002
$ws.Range("I81").Value = 506.5
$ws.Range("K81").Value = 1519.5
$ws.Range("M81").Value = -396.5

# CUL!row84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 10002
$ws.Range("I84").Value = 506.5
$ws.Range("K84").Value = 4558.5
$ws.Range("M84").Value = 1057.5

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4001624.8
$ws.Range("I131").Value = 704.75
$ws.Range("J131").Value = 7694782
$ws.Range("K131").Value = 2114.25
$ws.Range("L131").Value = 23084346
$ws.Range("M131").Value = 2925.75
$ws.Range("N131").Value = -23094426

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 71429550
$ws.Range("I132").Value = 250000450
$ws.Range("J132").Value = 1197.8
$ws.Range("K132").Value = 2250004050
$ws.Range("L132").Value = 10780.2
$ws.Range("M132").Value = -2250001520
$ws.Range("N132").Value = -15840.2

# GSM!row23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 837.5
$ws.Range("I23").Value = 828.5714
$ws.Range("K23").Value = 828.5714
$ws.Range("M23").Value = -605.5714

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2477.318
$ws.Range("I122").Value = 2254.8
$ws.Range("K122").Value = 6764.400000000001
$ws.Range("M122").Value = -4314.400000000001

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3118.4
$ws.Range("I132").Value = 2139.353
$ws.Range("K132").Value = 6418.059
$ws.Range("M132").Value = -3888.059

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 27224.75
$ws.Range("J46").Value = 99999
$ws.Range("L46").Value = 99999
$ws.Range("N46").Value = -100375

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8941.583000000001
$ws.Range("I136").Value = 7588.778
$ws.Range("J136").Value = 13000
$ws.Range("K136").Value = 22766.334
$ws.Range("L136").Value = 39000
$ws.Range("M136").Value = -20216.334
$ws.Range("N136").Value = -44100
